$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new blank column before the
# existing "Late" column (old column N), shifting Late/heading/Outstanding
# one column to the right (N->O, O->P, P->Q), and make this sheet active ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

$oldWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $oldWidth

$ws.Range("K14").Select() | Out-Null
